# uft-one-sap-fiori-general-ledger-overview/Default.xlsx
#
# Add a new "dtCard" column to the Global parameter sheet, holding the
# name of the Fiori app card/tile the script verifies
# ("Journal Entries to Be Verified"), and move the selection to the new
# cell below the table -- matching the state where the script finally
# runs without error.

$wb = $excel.ActiveWorkbook

$global = $wb.Worksheets.Item("Global")
$glOverview = $wb.Worksheets.Item("GLOverview")

$global.Activate()

# New header + value in column E, right after the existing
# BrowserName / URL / Username / Password columns.
$global.Range("E1").Value = "dtCard"
$global.Range("E2").Value = "Journal Entries to Be Verified"

# The last column of the 2-row parameter table carries the heavier
# "outer box" border (style of the old D2) while inner columns only
# get top/bottom rules (style of A2:C2). Re-use the existing cell
# formats instead of building new ones: copy the box border from D2
# onto the new last cell E2, then restore D2 to the plain inner style.
$global.Range("D2").Copy()
$global.Range("E2").PasteSpecial(-4122)

$global.Range("C2").Copy()
$global.Range("D2").PasteSpecial(-4122)

$global.Application.CutCopyMode = $false

$global.Columns.Item(5).ColumnWidth = 23.1

[void]$global.Range("E3").Select()

# Restore the workbook's original active sheet.
$glOverview.Activate()
